# Generate Report for Handback
# Update the "Latest Handoff Datetime" (column D) and "Latest Handback DateTime"
# (column G) for the first data row (the 36491f5a-... file) on both the
# "zh-cn" and "de-de" worksheets, reflecting a newer handoff/handback cycle.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D2").Value = "2016-03-09 03:13:30"
$ws_zhcn.Range("G2").Value = "2016-03-09 03:14:23"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D2").Value = "2016-03-09 03:13:41"
$ws_dede.Range("G2").Value = "2016-03-09 03:14:59"
